$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 39832
$ws.Range("D2").Value = 57545277
$ws.Range("C3").Value = 95204
$ws.Range("D3").Value = 139510276
$ws.Range("C4").Value = 32454
$ws.Range("D4").Value = 48045053
$ws.Range("C5").Value = 9172
$ws.Range("D5").Value = 13629008
$ws.Range("C6").Value = 2170
$ws.Range("D6").Value = 3223148
$ws.Range("C7").Value = 196
$ws.Range("D7").Value = 289093
$ws.Range("C12").Value = 43192
$ws.Range("D12").Value = 58532014
$ws.Range("C13").Value = 10132
$ws.Range("D13").Value = 14641111
$ws.Range("C14").Value = 26998
$ws.Range("D14").Value = 39569623
$ws.Range("C15").Value = 8603
$ws.Range("D15").Value = 12767484
$ws.Range("C16").Value = 2262
$ws.Range("D16").Value = 3360653
$ws.Range("C17").Value = 446
$ws.Range("D17").Value = 658123
$ws.Range("C20").Value = 10632
$ws.Range("D20").Value = 14039159
$ws.Range("C21").Value = 13991
$ws.Range("D21").Value = 20181922
$ws.Range("C22").Value = 32853
$ws.Range("D22").Value = 48185134
$ws.Range("C23").Value = 10583
$ws.Range("D23").Value = 15727732
$ws.Range("C24").Value = 2757
$ws.Range("D24").Value = 4100174
$ws.Range("C25").Value = 565
$ws.Range("D25").Value = 841592
$ws.Range("C27").Value = 12152
$ws.Range("D27").Value = 16191909
$ws.Range("C28").Value = 8087
$ws.Range("D28").Value = 11695234
$ws.Range("C29").Value = 23479
$ws.Range("D29").Value = 34452892
$ws.Range("C30").Value = 8082
$ws.Range("D30").Value = 12015633
$ws.Range("C31").Value = 2045
$ws.Range("D31").Value = 3050699
$ws.Range("C34").Value = 8681
$ws.Range("D34").Value = 11461889
$ws.Range("C35").Value = 3463
$ws.Range("D35").Value = 4998597
$ws.Range("C36").Value = 8244
$ws.Range("D36").Value = 12043599
$ws.Range("C37").Value = 3294
$ws.Range("D37").Value = 4883961
$ws.Range("C39").Value = 174
$ws.Range("D39").Value = 258686
$ws.Range("C41").Value = 2612
$ws.Range("D41").Value = 3530263
$ws.Range("C42").Value = 18128
$ws.Range("D42").Value = 26186802
$ws.Range("C43").Value = 53205
$ws.Range("D43").Value = 77961024
$ws.Range("C44").Value = 19576
$ws.Range("D44").Value = 29066453
$ws.Range("C45").Value = 5856
$ws.Range("D45").Value = 8715935
$ws.Range("C46").Value = 1310
$ws.Range("D46").Value = 1954644
$ws.Range("C50").Value = 17532
$ws.Range("D50").Value = 23257227
$ws.Range("C51").Value = 2232
$ws.Range("D51").Value = 3238783
$ws.Range("C52").Value = 7511
$ws.Range("D52").Value = 11035971
$ws.Range("C53").Value = 2499
$ws.Range("D53").Value = 3730184
$ws.Range("C54").Value = 787
$ws.Range("D54").Value = 1175415
$ws.Range("C57").Value = 7595
$ws.Range("D57").Value = 10446671
$ws.Range("C58").Value = 1358
$ws.Range("D58").Value = 2494131
$ws.Range("C59").Value = 3329
$ws.Range("D59").Value = 6128590
$ws.Range("C60").Value = 1316
$ws.Range("D60").Value = 2433866
$ws.Range("C61").Value = 443
$ws.Range("D61").Value = 816083
$ws.Range("C62").Value = 150
$ws.Range("D62").Value = 287600
$ws.Range("C64").Value = 2046
$ws.Range("D64").Value = 3467520
$ws.Range("C65").Value = 16245
$ws.Range("D65").Value = 23454417
$ws.Range("C66").Value = 46718
$ws.Range("D66").Value = 68315298
$ws.Range("C67").Value = 16304
$ws.Range("D67").Value = 24223549
$ws.Range("C68").Value = 4758
$ws.Range("D68").Value = 7086788
$ws.Range("C69").Value = 1013
$ws.Range("D69").Value = 1505829
$ws.Range("C70").Value = 85
$ws.Range("D70").Value = 124830
$ws.Range("C73").Value = 15683
$ws.Range("D73").Value = 20619162
$ws.Range("C74").Value = 57433
$ws.Range("D74").Value = 83521746
$ws.Range("C75").Value = 158993
$ws.Range("D75").Value = 234091497
$ws.Range("C76").Value = 68186
$ws.Range("D76").Value = 101571719
$ws.Range("C77").Value = 21973
$ws.Range("D77").Value = 32832013
$ws.Range("C78").Value = 5328
$ws.Range("D78").Value = 7958102
$ws.Range("C79").Value = 330
$ws.Range("D79").Value = 490170
$ws.Range("C85").Value = 56445
$ws.Range("D85").Value = 76438914
$ws.Range("C86").Value = 4916
$ws.Range("D86").Value = 7123839
$ws.Range("C87").Value = 12187
$ws.Range("D87").Value = 17900881
$ws.Range("C88").Value = 4031
$ws.Range("D88").Value = 6006458
$ws.Range("C89").Value = 1398
$ws.Range("D89").Value = 2088111
$ws.Range("C90").Value = 312
$ws.Range("D90").Value = 465512
$ws.Range("C93").Value = 5703
$ws.Range("D93").Value = 7660268
$ws.Range("C94").Value = 1716
$ws.Range("D94").Value = 2473136
$ws.Range("C95").Value = 5539
$ws.Range("D95").Value = 8160504
$ws.Range("C96").Value = 2029
$ws.Range("D96").Value = 3020308
$ws.Range("C97").Value = 739
$ws.Range("D97").Value = 1107460
$ws.Range("C101").Value = 3791
$ws.Range("D101").Value = 5025440
$ws.Range("C102").Value = 796
$ws.Range("D102").Value = 1422215
$ws.Range("C103").Value = 505
$ws.Range("D103").Value = 939892
$ws.Range("C107").Value = 11379
$ws.Range("D107").Value = 16498883
$ws.Range("C108").Value = 30288
$ws.Range("D108").Value = 44472764
$ws.Range("C109").Value = 10153
$ws.Range("D109").Value = 15093945
$ws.Range("C110").Value = 2802
$ws.Range("D110").Value = 4177510
$ws.Range("C111").Value = 526
$ws.Range("D111").Value = 783403
$ws.Range("C114").Value = 10189
$ws.Range("D114").Value = 13430504
$ws.Range("C115").Value = 32020
$ws.Range("D115").Value = 46149833
$ws.Range("C116").Value = 68792
$ws.Range("D116").Value = 100638176
$ws.Range("C117").Value = 22102
$ws.Range("D117").Value = 32836804
$ws.Range("C118").Value = 6296
$ws.Range("D118").Value = 9374011
$ws.Range("C119").Value = 1200
$ws.Range("D119").Value = 1792889
$ws.Range("C120").Value = 104
$ws.Range("D120").Value = 152395
$ws.Range("C124").Value = 26848
$ws.Range("D124").Value = 35809087
$ws.Range("C125").Value = 38058
$ws.Range("D125").Value = 54889371
$ws.Range("C126").Value = 80347
$ws.Range("D126").Value = 117447950
$ws.Range("C127").Value = 24797
$ws.Range("D127").Value = 36801270
$ws.Range("C128").Value = 6672
$ws.Range("D128").Value = 9915243
$ws.Range("C129").Value = 1355
$ws.Range("D129").Value = 2010958
$ws.Range("C133").Value = 33166
$ws.Range("D133").Value = 43981558
$ws.Range("C134").Value = 13983
$ws.Range("D134").Value = 20235093
$ws.Range("C135").Value = 33643
$ws.Range("D135").Value = 49395272
$ws.Range("C136").Value = 11879
$ws.Range("D136").Value = 17648784
$ws.Range("C137").Value = 3117
$ws.Range("D137").Value = 4646741
$ws.Range("C138").Value = 535
$ws.Range("D138").Value = 796490
$ws.Range("C141").Value = 11249
$ws.Range("D141").Value = 14976934
$ws.Range("C142").Value = 37155
$ws.Range("D142").Value = 53655018
$ws.Range("C143").Value = 85461
$ws.Range("D143").Value = 125172678
$ws.Range("C144").Value = 25426
$ws.Range("D144").Value = 37773178
$ws.Range("C145").Value = 6681
$ws.Range("D145").Value = 9966025
$ws.Range("C146").Value = 1533
$ws.Range("D146").Value = 2278802
$ws.Range("C149").Value = 30552
$ws.Range("D149").Value = 41134172
